# Applies the numeric value updates to the Siren_Profits workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
# Each worksheet backs an Excel Table (Table_<SHEET>) with columns A:N; only data cells in
# columns H-N (price/profit figures) are updated here, matching the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 565.7
$ws.Range("I2").Value = 466.14285
$ws.Range("K2").Value = 466.14285
$ws.Range("M2").Value = -353.14285
$ws.Range("H17").Value = 1368.7778
$ws.Range("I17").Value = 1299
$ws.Range("K17").Value = 3897
$ws.Range("M17").Value = -3729
$ws.Range("J33").Value = 1510
$ws.Range("I33").Value = 299.3
$ws.Range("L33").Value = 1510
$ws.Range("K33").Value = 299.3
$ws.Range("H33").Value = 753.3125
$ws.Range("M33").Value = -70.30000000000001
$ws.Range("N33").Value = -1968
$ws.Range("J46").Value = 1000000
$ws.Range("H46").Value = 1000000
$ws.Range("N46").Value = -3000238
$ws.Range("L46").Value = 3000000
$ws.Range("H60").Value = 1000000
$ws.Range("L60").Value = 3000000
$ws.Range("N60").Value = -3000968
$ws.Range("J60").Value = 1000000
$ws.Range("H74").Value = 3240.1177
$ws.Range("I74").Value = 2386.375
$ws.Range("L74").Value = 3999
$ws.Range("K74").Value = 2386.375
$ws.Range("J74").Value = 3999
$ws.Range("M74").Value = -1450.375
$ws.Range("N74").Value = -5871
$ws.Range("K76").Value = 3754
$ws.Range("M76").Value = -3439
$ws.Range("H76").Value = 4176
$ws.Range("I76").Value = 3754
$ws.Range("M77").Value = -7251.875
$ws.Range("I77").Value = 2386.375
$ws.Range("L77").Value = 19995
$ws.Range("H77").Value = 3240.1177
$ws.Range("N77").Value = -29355
$ws.Range("K77").Value = 11931.875
$ws.Range("J77").Value = 3999
$ws.Range("I79").Value = 3754
$ws.Range("H79").Value = 4176
$ws.Range("K79").Value = 3754
$ws.Range("M79").Value = -2662
$ws.Range("N86").Value = -147446
$ws.Range("L86").Value = 145200
$ws.Range("H86").Value = 56354.875
$ws.Range("J86").Value = 145200
$ws.Range("L89").Value = 726000
$ws.Range("N89").Value = -737232
$ws.Range("J89").Value = 145200
$ws.Range("H89").Value = 56354.875
$ws.Range("I100").Value = 117396.16
$ws.Range("H100").Value = 96863.09
$ws.Range("M100").Value = -116855.16
$ws.Range("K100").Value = 117396.16
$ws.Range("K107").Value = 9013.111000000001
$ws.Range("I107").Value = 9013.111000000001
$ws.Range("H107").Value = 6579.1875
$ws.Range("M107").Value = -7093.111000000001
$ws.Range("L112").Value = 7270.3842
$ws.Range("J112").Value = 2423.4614
$ws.Range("M112").Value = -3586.5002
$ws.Range("K112").Value = 4694.5002
$ws.Range("I112").Value = 1564.8334
$ws.Range("H112").Value = 2152.3157
$ws.Range("N112").Value = -9486.3842
$ws.Range("H113").Value = 13983
$ws.Range("J113").Value = 12470.75
$ws.Range("L113").Value = 12470.75
$ws.Range("N113").Value = -18978.75
$ws.Range("N118").Value = -6864.799999999999
$ws.Range("L118").Value = 3550.8
$ws.Range("J118").Value = 1183.6
$ws.Range("H118").Value = 911.875
$ws.Range("I137").Value = 14395
$ws.Range("J137").Value = 2075.2
$ws.Range("H137").Value = 8433.807000000001
$ws.Range("K137").Value = 43185
$ws.Range("M137").Value = -40635
$ws.Range("L137").Value = 6225.599999999999
$ws.Range("N137").Value = -11325.6
$ws.Range("N138").Value = -22147.5314
$ws.Range("L138").Value = 11867.5314
$ws.Range("J138").Value = 3955.8438
$ws.Range("H138").Value = 221966.86

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 92681.91
$ws.Range("I2").Value = 646.5714
$ws.Range("K2").Value = 646.5714
$ws.Range("M2").Value = -533.5714
$ws.Range("K21").Value = 9999.799999999999
$ws.Range("L21").Value = 5000
$ws.Range("J21").Value = 5000
$ws.Range("M21").Value = -9625.799999999999
$ws.Range("I21").Value = 9999.799999999999
$ws.Range("H21").Value = 9166.5
$ws.Range("N21").Value = -5748
$ws.Range("M45").Value = -121371
$ws.Range("L45").Value = 5490
$ws.Range("N45").Value = -6244
$ws.Range("J45").Value = 5490
$ws.Range("H45").Value = 78689.484
$ws.Range("K45").Value = 121748
$ws.Range("I45").Value = 121748
$ws.Range("H74").Value = 1481.8928
$ws.Range("I74").Value = 1291.72
$ws.Range("K74").Value = 1291.72
$ws.Range("M74").Value = -417.72
$ws.Range("M77").Value = -2090.6
$ws.Range("I77").Value = 1291.72
$ws.Range("H77").Value = 1481.8928
$ws.Range("K77").Value = 6458.6
$ws.Range("M116").Value = 1647.4286
$ws.Range("H116").Value = 92681.91
$ws.Range("K116").Value = 646.5714
$ws.Range("I116").Value = 646.5714
$ws.Range("H132").Value = 2277.016
$ws.Range("I132").Value = 1934.814
$ws.Range("K132").Value = 5804.442
$ws.Range("M132").Value = -3274.442
$ws.Range("J135").Value = 79980
$ws.Range("L135").Value = 79980
$ws.Range("H135").Value = 79980
$ws.Range("N135").Value = -90120

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 92681.91
$ws.Range("K3").Value = 646.5714
$ws.Range("I3").Value = 646.5714
$ws.Range("M3").Value = -532.5714
$ws.Range("K82").Value = 7399.5
$ws.Range("M82").Value = -7016.5
$ws.Range("I82").Value = 7399.5
$ws.Range("H82").Value = 20549.375
$ws.Range("N82").Value = -60765
$ws.Range("L82").Value = 59999
$ws.Range("J82").Value = 59999
$ws.Range("L85").Value = 59999
$ws.Range("M85").Value = -6073.5
$ws.Range("H85").Value = 20549.375
$ws.Range("K85").Value = 7399.5
$ws.Range("I85").Value = 7399.5
$ws.Range("N85").Value = -62651
$ws.Range("J85").Value = 59999
$ws.Range("N86").Value = -4559
$ws.Range("K86").Value = 7253
$ws.Range("L86").Value = 2313
$ws.Range("I86").Value = 7253
$ws.Range("H86").Value = 5953
$ws.Range("J86").Value = 2313
$ws.Range("M86").Value = -6130
$ws.Range("L89").Value = 11565
$ws.Range("M89").Value = -30649
$ws.Range("K89").Value = 36265
$ws.Range("N89").Value = -22797
$ws.Range("I89").Value = 7253
$ws.Range("J89").Value = 2313
$ws.Range("H89").Value = 5953
$ws.Range("H94").Value = 11519.429
$ws.Range("N94").Value = -5860.125
$ws.Range("K94").Value = 14143.95
$ws.Range("L94").Value = 4958.125
$ws.Range("M94").Value = -13692.95
$ws.Range("I94").Value = 14143.95
$ws.Range("J94").Value = 4958.125
$ws.Range("H99").Value = 16400.217
$ws.Range("K99").Value = 20381.25
$ws.Range("I99").Value = 20381.25
$ws.Range("J99").Value = 7300.7144
$ws.Range("M99").Value = -18883.25
$ws.Range("L99").Value = 7300.7144
$ws.Range("N99").Value = -10296.7144
$ws.Range("I134").Value = 5560.4644
$ws.Range("H134").Value = 4927.0884
$ws.Range("M134").Value = -14146.3932
$ws.Range("K134").Value = 16681.3932

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -22969
$ws.Range("L31").Value = 5794.8
$ws.Range("H31").Value = 15323.454
$ws.Range("J31").Value = 5794.8
$ws.Range("K31").Value = 23264
$ws.Range("N31").Value = -6384.8
$ws.Range("I31").Value = 23264
$ws.Range("I34").Value = 23264
$ws.Range("J34").Value = 5794.8
$ws.Range("L34").Value = 5794.8
$ws.Range("K34").Value = 23264
$ws.Range("N34").Value = -6198.8
$ws.Range("H34").Value = 15323.454
$ws.Range("M34").Value = -23062
$ws.Range("K58").Value = 2412.9583
$ws.Range("I58").Value = 2412.9583
$ws.Range("M58").Value = -2209.9583
$ws.Range("H58").Value = 2418.4443
$ws.Range("H105").Value = 236917.44
$ws.Range("I105").Value = 303436.72
$ws.Range("M105").Value = -301689.72
$ws.Range("K105").Value = 303436.72
$ws.Range("H132").Value = 1848.6818
$ws.Range("I132").Value = 1888.55
$ws.Range("K132").Value = 5665.65
$ws.Range("M132").Value = -3135.65
$ws.Range("L132").Value = 4350
$ws.Range("N132").Value = -9410
$ws.Range("J132").Value = 1450
$ws.Range("L134").Value = 5620.5
$ws.Range("I134").Value = 2361.5925
$ws.Range("J134").Value = 1873.5
$ws.Range("N134").Value = -10690.5
$ws.Range("H134").Value = 2229.6758
$ws.Range("M134").Value = -4549.7775
$ws.Range("K134").Value = 7084.7775
$ws.Range("K136").Value = 7238.874899999999
$ws.Range("M136").Value = -4688.874899999999
$ws.Range("H136").Value = 2418.4443
$ws.Range("I136").Value = 2412.9583
$ws.Range("H141").Value = 347419
$ws.Range("J141").Value = 471827.78
$ws.Range("N141").Value = -482187.78
$ws.Range("L141").Value = 471827.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J7").Value = 112.25
$ws.Range("L7").Value = 336.75
$ws.Range("I7").Value = 25.375
$ws.Range("N7").Value = -560.75
$ws.Range("H7").Value = 54.333332
$ws.Range("M7").Value = 35.875
$ws.Range("K7").Value = 76.125
$ws.Range("K35").Value = 2380.125
$ws.Range("I35").Value = 793.375
$ws.Range("M35").Value = -2092.125
$ws.Range("H35").Value = 890.63635
$ws.Range("M56").Value = -5654.4116
$ws.Range("I56").Value = 6184.4116
$ws.Range("H56").Value = 6184.4116
$ws.Range("K56").Value = 6184.4116
$ws.Range("L107").Value = 4834.3845
$ws.Range("N107").Value = -8674.3845
$ws.Range("K107").Value = 1162.5
$ws.Range("J107").Value = 1611.4615
$ws.Range("I107").Value = 387.5
$ws.Range("H107").Value = 1145.1904
$ws.Range("M107").Value = 757.5
$ws.Range("M117").Value = -10550
$ws.Range("L117").Value = 15000
$ws.Range("N117").Value = -21884
$ws.Range("J117").Value = 5000
$ws.Range("H117").Value = 4748
$ws.Range("I117").Value = 4664
$ws.Range("K117").Value = 13992
$ws.Range("J129").Value = 2622.9285
$ws.Range("K129").Value = 1896.6
$ws.Range("I129").Value = 632.2
$ws.Range("M129").Value = 3103.4
$ws.Range("H129").Value = 1793.4584
$ws.Range("N129").Value = -17868.7855
$ws.Range("L129").Value = 7868.7855
$ws.Range("H139").Value = 1252557.5
$ws.Range("L139").Value = 12000
$ws.Range("J139").Value = 4000
$ws.Range("N139").Value = -22280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N80").Value = -9055.666499999999
$ws.Range("J80").Value = 7059.6665
$ws.Range("H80").Value = 8988.723
$ws.Range("L80").Value = 7059.6665
$ws.Range("H83").Value = 8988.723
$ws.Range("L83").Value = 35298.3325
$ws.Range("N83").Value = -45282.3325
$ws.Range("J83").Value = 7059.6665
$ws.Range("J110").Value = 90555
$ws.Range("H110").Value = 90555
$ws.Range("N110").Value = -98735
$ws.Range("L110").Value = 90555
$ws.Range("H113").Value = 29495.25
$ws.Range("M113").Value = -35823.668
$ws.Range("J113").Value = 4000
$ws.Range("I113").Value = 37993.668
$ws.Range("L113").Value = 4000
$ws.Range("K113").Value = 37993.668
$ws.Range("N113").Value = -8340
$ws.Range("K122").Value = 17109.4731
$ws.Range("H122").Value = 8747.037
$ws.Range("M122").Value = -14659.4731
$ws.Range("I122").Value = 5703.1577
$ws.Range("H132").Value = 2250.6418
$ws.Range("I132").Value = 2159.0544
$ws.Range("K132").Value = 6477.1632
$ws.Range("M132").Value = -3947.1632

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J7").Value = 7332.778
$ws.Range("L7").Value = 7332.778
$ws.Range("I7").Value = 30334.941
$ws.Range("N7").Value = -7556.778
$ws.Range("H7").Value = 22372.654
$ws.Range("M7").Value = -30222.941
$ws.Range("K7").Value = 30334.941
$ws.Range("I40").Value = 27503.588
$ws.Range("K40").Value = 27503.588
$ws.Range("M40").Value = -27367.588
$ws.Range("H40").Value = 22028.867
$ws.Range("J46").Value = 5365233.5
$ws.Range("I46").Value = 800
$ws.Range("K46").Value = 800
$ws.Range("H46").Value = 3671201.8
$ws.Range("N46").Value = -5365609.5
$ws.Range("L46").Value = 5365233.5
$ws.Range("M46").Value = -612
$ws.Range("K82").Value = 3774
$ws.Range("M82").Value = -3413
$ws.Range("I82").Value = 3774
$ws.Range("H82").Value = 2665.4
$ws.Range("N82").Value = -2417.375
$ws.Range("L82").Value = 1695.375
$ws.Range("J82").Value = 1695.375
$ws.Range("L85").Value = 1695.375
$ws.Range("M85").Value = -2526
$ws.Range("H85").Value = 2665.4
$ws.Range("K85").Value = 3774
$ws.Range("I85").Value = 3774
$ws.Range("N85").Value = -4191.375
$ws.Range("J85").Value = 1695.375
$ws.Range("K93").Value = 6568.9287
$ws.Range("I93").Value = 6568.9287
$ws.Range("H93").Value = 5312.5264
$ws.Range("M93").Value = -5320.9287
$ws.Range("H126").Value = 22372.654
$ws.Range("N126").Value = -26938.334
$ws.Range("L126").Value = 21998.334
$ws.Range("I126").Value = 30334.941
$ws.Range("M126").Value = -88534.823
$ws.Range("J126").Value = 7332.778
$ws.Range("K126").Value = 91004.823
$ws.Range("H132").Value = 830003.25
$ws.Range("I132").Value = 933144.25
$ws.Range("K132").Value = 2799432.75
$ws.Range("M132").Value = -2796902.75
$ws.Range("L132").Value = 14625
$ws.Range("N132").Value = -19685
$ws.Range("J132").Value = 4875
$ws.Range("J135").Value = 97266.53999999999
$ws.Range("L135").Value = 97266.53999999999
$ws.Range("H135").Value = 97266.53999999999
$ws.Range("N135").Value = -107406.54
$ws.Range("L136").Value = 15893.4552
$ws.Range("J136").Value = 5297.8184
$ws.Range("K136").Value = 7684.399800000001
$ws.Range("M136").Value = -5134.399800000001
$ws.Range("H136").Value = 3719.1538
$ws.Range("I136").Value = 2561.4666
$ws.Range("N136").Value = -20993.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K81").Value = 26996
$ws.Range("N81").Value = -8415.6666
$ws.Range("L81").Value = 6293.6666
$ws.Range("I81").Value = 13498
$ws.Range("M81").Value = -25935
$ws.Range("J81").Value = 3146.8333
$ws.Range("H81").Value = 9061.786
$ws.Range("J84").Value = 3146.8333
$ws.Range("H84").Value = 9061.786
$ws.Range("I84").Value = 13498
$ws.Range("L84").Value = 31468.333
$ws.Range("K84").Value = 134980
$ws.Range("N84").Value = -42076.333
$ws.Range("M84").Value = -129676
$ws.Range("L112").Value = 46650
$ws.Range("J112").Value = 46650
$ws.Range("H112").Value = 46650
$ws.Range("N112").Value = -49604
$ws.Range("K122").Value = 7818
$ws.Range("H122").Value = 5548.275
$ws.Range("M122").Value = -5368
$ws.Range("I122").Value = 2606
$ws.Range("H132").Value = 10004.02
$ws.Range("I132").Value = 13231.147
$ws.Range("K132").Value = 39693.44100000001
$ws.Range("M132").Value = -37163.44100000001
$ws.Range("L132").Value = 9439.125
$ws.Range("N132").Value = -14499.125
$ws.Range("J132").Value = 3146.375
$ws.Range("L136").Value = 17151.375
$ws.Range("J136").Value = 5717.125
$ws.Range("K136").Value = 1327580.34
$ws.Range("M136").Value = -1325030.34
$ws.Range("H136").Value = 361259.88
$ws.Range("I136").Value = 442526.78
$ws.Range("N136").Value = -22251.375
